$d = $word.ActiveDocument

# --- Change 1: "Objetivos" paragraph -------------------------------------
# "con el fin de implementar dicha plataforma" -> "con el fin de utilizar dicha plataforma"
$d.Content.Find.Execute(
    "con el fin de implementar dicha plataforma",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "con el fin de utilizar dicha plataforma", 2) | Out-Null

# --- Change 2: "Conclusiones" paragraph rewrite ---------------------------
# Replace the opening clause.
$d.Content.Find.Execute(
    "Esta investigación se contuará en proceso todavía, pero podemos manifestar con los aspectos abarcados hasta la actualidad que",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Podemos afirmar  que", 2) | Out-Null

# Rework the "juego" sentence.
$d.Content.Find.Execute(
    "Si es juego deberá se dinámica y activa todo el tiempo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Si hablamos de un juego deberá ser dinámica y activa todo el tiempo, generando un enfoque absoluto de los sentidos del usuario.",
    2) | Out-Null

# Fix the typo "usaurio" -> "usuario".
$d.Content.Find.Execute(
    "comodidad del usaurio.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "comodidad del usuario.", 2) | Out-Null

# --- Append a new closing paragraph, then a trailing blank paragraph -----
$rng = $d.Content
$rng.Find.Execute("comodidad del usuario.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraRange = $rng.Paragraphs(1).Range
$paraRange.InsertParagraphAfter()

$newPara = $d.Range($paraRange.End, $paraRange.End)
$newPara.Text = "Muchas de las actividades que posteriormente puedan generarse con el uso de HCI deberán tener en cuenta este tipo de análisis, siendo que pueden aplicar cambios estructurales e innovación a las formas de interacción."

$newPara.InsertParagraphAfter()
